# Apply FHIR IG terminology and profile corrections to the Metadata sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Metadata")

# Experimental row (row 7): value column was blank -> literal text "false".
# A plain .Value assignment of "false" gets auto-coerced to a Boolean by
# Excel's type inference, so stage the literal text in a scratch cell
# (formatted as Text) and paste-special the *value* into B7 -- this keeps
# B7 a genuine text cell (and preserves its existing style) instead of a
# Boolean cell.
$scratch = $ws.Range("Z1")
$scratch.NumberFormat = "@"
$scratch.Value = "'false"
$scratch.Copy() | Out-Null
$ws.Range("B7").PasteSpecial(-4163) | Out-Null
$scratch.Clear() | Out-Null

# Date row (row 8): refresh the generation timestamp.
$ws.Range("B8").Value = "2025-11-30T13:08:37+00:00"

# Description row (row 17): add the value set description.
$ws.Range("B17").Value = "Categories for interpreting recovery readiness scores from wearable devices"
